$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 header: end date moves to 2020-08-07 and the staff name is cleared
# (task list now shows the "Individual" column instead of a single staff name).
$ws.Range("A2").Value = "Start Date: 2020-08-05, End Date: 2020-08-07, Staff name: "

# Row 5: now reflects staff id 2 / "Kidden", entry dated 2020-08-07, entry No. 2,
# with updated calendar/adjust/work hours.
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Kidden"

# C5 must stay literal text "2020-08-07" (not get auto-parsed into a date
# serial number), while keeping the default/unstyled cell format - force the
# cell to Text first, write the value, then restore the original (default)
# style so no new style entry is introduced.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2020-08-07"
$ws.Range("C5").Style = $ws.Range("B5").Style

$ws.Range("D5").Value = "entryname No. 2 on "
$ws.Range("E5").Value = "01 : 41"
$ws.Range("F5").Value = "00 : 40"
$ws.Range("G5").Value = 2.35
